$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sharedStrings text changes ---
# Row 30's rule_description (also used for row numbers with threshold referencing index 63)
$ws.Range("D30").Value = "swh_max_swan > t1 AND anom_swh_p80_waverys > t2"

# Row 30's threshold text (np.float64 tuple) index 74
$ws.Range("E30").Value = "(np.float64(1.9166666), np.float64(0.1887734374999999))"

# --- Row 2 updates ---
$ws.Range("N2").Value = 797999.3275
$ws.Range("O2").Value = 1008950

# --- Row 21 updates ---
$ws.Range("O21").Value = 4263519.857142854
$ws.Range("P21").Value = 102.51975
$ws.Range("Q21").Value = 135
$ws.Range("T21").Value = 83.944
$ws.Range("U21").Value = 114
$ws.Range("V21").Value = 18.57575
$ws.Range("W21").Value = 35
$ws.Range("X21").Value = 19.37175
$ws.Range("Y21").Value = 38

# --- Row 24 updates ---
$ws.Range("O24").Value = 322602.8571428572
$ws.Range("P24").Value = 89.682
$ws.Range("Q24").Value = 122
$ws.Range("T24").Value = 57.20575
$ws.Range("U24").Value = 81
$ws.Range("V24").Value = 32.47625
$ws.Range("W24").Value = 53
$ws.Range("X24").Value = 17.10175
$ws.Range("Y24").Value = 32

# --- Row 27 updates ---
$ws.Range("O27").Value = 369686.1428571429
$ws.Range("P27").Value = 121.133
$ws.Range("Q27").Value = 158
$ws.Range("T27").Value = 54.771
$ws.Range("U27").Value = 77
$ws.Range("V27").Value = 66.36199999999999
$ws.Range("W27").Value = 96
$ws.Range("X27").Value = 22.11075
$ws.Range("Y27").Value = 38

# --- Row 30 updates ---
$ws.Range("F30").Value = 406
$ws.Range("G30").Value = 330
$ws.Range("H30").Value = 1724
$ws.Range("I30").Value = 101
$ws.Range("J30").Value = 0.5516304347826086
$ws.Range("K30").Value = 0.8007889546351085
$ws.Range("L30").Value = 0.8317063647012886
$ws.Range("M30").Value = 0.6532582461786002
$ws.Range("N30").Value = 119993.4537142857
$ws.Range("O30").Value = 158719.9999999999
$ws.Range("P30").Value = 105.841
$ws.Range("Q30").Value = 140
$ws.Range("T30").Value = 47.75075
$ws.Range("U30").Value = 67
$ws.Range("V30").Value = 58.09025
$ws.Range("W30").Value = 85
$ws.Range("X30").Value = 14.265
$ws.Range("Y30").Value = 25
